# Atualização das Planilhas de Execução
# Fill in "Tempo Gasto/min" (column L) with 10 minutes for the rows
# that were worked on in this update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(8,9,10,11,12,13,14,15,16,17,18,19,20,21,27,28,29,30,54,55,56,57,58,59,60,61,62,63,64,65,66,67)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 12).Value = 10
}

# Reflect the cursor/selection position left in the sheet after the edits.
$ws.Range("L50").Select()
